$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J, formatted like the existing header cells (copy style from H1).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data rows 2-68: plain numeric values in columns I (I0) and J (IF).
$data = @{
    2 = @(8, 8)
    3 = @(8, 8)
    4 = @(8, 8)
    5 = @(6, 6)
    6 = @(8, 9)
    7 = @(9, 9)
    8 = @(5, 6)
    9 = @(7, 8)
    10 = @(5, 6)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(6, 6)
    14 = @(5, 6)
    15 = @(1, 2)
    16 = @(6, 7)
    17 = @(6, 7)
    18 = @(7, 7)
    19 = @(9, 9)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(7, 7)
    23 = @(5, 6)
    24 = @(9, 9)
    25 = @(7, 7)
    26 = @(8, 8)
    27 = @(11, 12)
    28 = @(8, 8)
    29 = @(8, 8)
    30 = @(8, 8)
    31 = @(7, 7)
    32 = @(6, 6)
    33 = @(7, 8)
    34 = @(9, 9)
    35 = @(8, 8)
    36 = @(8, 8)
    37 = @(8, 8)
    38 = @(5, 5)
    39 = @(8, 8)
    40 = @(8, 8)
    41 = @(7, 8)
    42 = @(6, 7)
    43 = @(6, 7)
    44 = @(8, 8)
    45 = @(6, 6)
    46 = @(10, 11)
    47 = @(7, 7)
    48 = @(5, 5)
    49 = @(9, 9)
    50 = @(6, 6)
    51 = @(7, 7)
    52 = @(10, 10)
    53 = @(8, 8)
    54 = @(8, 8)
    55 = @(7, 7)
    56 = @(7, 7)
    57 = @(6, 6)
    58 = @(8, 8)
    59 = @(9, 9)
    60 = @(6, 6)
    61 = @(6, 6)
    62 = @(6, 6)
    63 = @(6, 7)
    64 = @(5, 5)
    65 = @(1, 1)
    66 = @(7, 8)
    67 = @(6, 6)
    68 = @(3, 3)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}

Write-Output "done"
